$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "44.319.96"
$ws.Range("E2").Value = "  -0.01%  "
$ws.Range("D3").Value = "2.265.60"
$ws.Range("E3").Value = "  -0.27%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'316.67"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.10%  "
$ws.Range("D6").Value = "'100.01"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -5.01%  "
$ws.Range("E7").Value = "  -2.17%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("E9").Value = "  -5.30%  "
$ws.Range("D10").Value = "'36.38"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -6.23%  "
$ws.Range("E11").Value = "  -1.73%  "
$ws.Range("D12").Value = "'7.45"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -5.61%  "
$ws.Range("E13").Value = "  -2.54%  "
$ws.Range("D14").Value = "2.607.09"
$ws.Range("E14").Value = "  -0.45%  "
$ws.Range("D15").Value = "'0.852"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.70%  "
$ws.Range("D16").Value = "2.257.03"
$ws.Range("E16").Value = "  -0.73%  "
$ws.Range("E17").Value = "  -3.57%  "
$ws.Range("D18").Value = "44.202.06"
$ws.Range("E18").Value = "  -0.10%  "
$ws.Range("D19").Value = "'13.27"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -4.25%  "
$ws.Range("D20").Value = "0.0₃0984"
$ws.Range("E20").Value = "  -2.03%  "
$ws.Range("D21").Value = "'6.43"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.58%  "
$ws.Range("D22").Value = "'65.92"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Value = "'239.76"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.41%  "
$ws.Range("D24").Value = "'3.01"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -6.27%  "
$ws.Range("D25").Value = "'2.06"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -7.38%  "
$ws.Range("D26").Value = "'1.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.44%  "
$ws.Range("D27").Value = "'10.33"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.54%  "
$ws.Range("D28").Value = "'38.66"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.65%  "
$ws.Range("D29").Value = "'2.13"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.75%  "
$ws.Range("D30").Value = "'6.14"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -5.57%  "
$ws.Range("D31").Value = "'20.22"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.21%  "
$ws.Range("E32").Value = "  +12.96%  "
$ws.Range("B33").Value = "Monero"
$ws.Range("C33").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D33").Value = "'155.01"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -5.53%  "
$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").Value = "'0.0847"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.99%  "
$ws.Range("E35").Value = "  -3.19%  "
$ws.Range("D36").Value = "'1.94"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.87%  "
$ws.Range("E37").Value = "  -5.42%  "
$ws.Range("E38").Value = "  -2.18%  "
$ws.Range("D39").Value = "'15.53"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.39%  "
$ws.Range("E40").Value = "  -9.91%  "
$ws.Range("E41").Value = "  -9.17%  "
$ws.Range("E42").Value = "  -5.48%  "
$ws.Range("E43").Value = "  +0.10%  "
$ws.Range("D44").Value = "1.733.85"
$ws.Range("E44").Value = "  -2.56%  "
$ws.Range("D45").Value = "'83.97"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.61%  "
$ws.Range("D46").Value = "'0.198"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -5.31%  "
$ws.Range("E47").Value = "  -3.68%  "
$ws.Range("B48").Value = "Aave"
$ws.Range("C48").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D48").Value = "'102.85"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.27%  "
$ws.Range("B49").Value = "ordi"
$ws.Range("C49").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D49").Value = "'72.38"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.16%  "
$ws.Range("D50").Value = "'57.30"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -5.59%  "
$ws.Range("B51").Value = "FraxShare"
$ws.Range("C51").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D51").Value = "'8.33"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.36%  "
